# Updates cryptos list price/volume data (scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the literal text into the cell without Excel re-interpreting
    # numeric-looking strings (e.g. "214.50") as numbers, while leaving the
    # cells applied style/format untouched (matches source inlineStr cells).
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '25.908.14'
$ws.Range('E2').Value = '  +0.04%  '

$ws.Range('D3').Value = '1.636.51'
$ws.Range('E3').Value = '  +0.11%  '

$ws.Range('E4').Value = '  +0.25%  '

Set-TextValue 'D5' '214.50'
$ws.Range('E5').Value = '  -0.18%  '

$ws.Range('E6').Value = '  +1.07%  '

$ws.Range('E7').Value = '  +0.27%  '

$ws.Range('E8').Value = '  -0.87%  '

$ws.Range('E9').Value = '  +0.44%  '

Set-TextValue 'D10' '19.59'
$ws.Range('E10').Value = '  -0.44%  '

$ws.Range('E11').Value = '  +0.51%  '

$ws.Range('D12').Value = '1.863.26'
$ws.Range('E12').Value = '  +0.06%  '

$ws.Range('D14').Value = '1.653.32'
$ws.Range('E14').Value = '  -0.17%  '

$ws.Range('E15').Value = '  -1.71%  '

$ws.Range('E16').Value = '  -0.53%  '

Set-TextValue 'D17' '62.59'
$ws.Range('E17').Value = '  -0.49%  '

$ws.Range('D18').Value = '25.925.02'
$ws.Range('E18').Value = '  +0.12%  '

$ws.Range('E19').Value = '  +0.19%  '

Set-TextValue 'D20' '193.69'
$ws.Range('E20').Value = '  +1.05%  '

$ws.Range('E21').Value = '  -1.23%  '

$ws.Range('E22').Value = '  -0.57%  '

$ws.Range('E23').Value = '  -0.95%  '

$ws.Range('E24').Value = '  +0.32%  '

Set-TextValue 'D25' '143.63'
$ws.Range('E25').Value = '  +0.94%  '

$ws.Range('E26').Value = '  +0.27%  '

$ws.Range('E27').Value = '  +2.70%  '

Set-TextValue 'D28' '6.83'
$ws.Range('E28').Value = '  -0.42%  '

$ws.Range('E29').Value = '  -0.66%  '

$ws.Range('E30').Value = '  +0.16%  '

$ws.Range('E31').Value = '  +1.16%  '

$ws.Range('E32').Value = '  -1.19%  '

$ws.Range('E33').Value = '  -1.00%  '

$ws.Range('E34').Value = '  -2.58%  '

$ws.Range('E35').Value = '  +1.43%  '

$ws.Range('E36').Value = '  -0.60%  '

$ws.Range('D37').Value = '1.138.75'
$ws.Range('E37').Value = '  -0.80%  '

Set-TextValue 'D38' '0.544'
$ws.Range('E38').Value = '  +0.03%  '

Set-TextValue 'D39' '2.48'
$ws.Range('E39').Value = '  -1.07%  '

$ws.Range('E40').Value = '  +0.05%  '

$ws.Range('E41').Value = '  +0.24%  '

Set-TextValue 'D42' '99.32'
$ws.Range('E42').Value = '  -1.41%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D43' '0.798'
$ws.Range('E43').Value = '  -0.59%  '

$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '5.42'
$ws.Range('E44').Value = '  -3.83%  '

$ws.Range('D45').Value = '1.772.06'
$ws.Range('E45').Value = '  -0.01%  '

$ws.Range('D46').Value = '0.0₆0115'
$ws.Range('E46').Value = '  +3.68%  '

$ws.Range('E47').Value = '  +1.37%  '

$ws.Range('E50').Value = '  -0.37%  '

Set-TextValue 'D51' '7.64'
$ws.Range('E51').Value = '  +0.11%  '
